$d = $word.ActiveDocument
$wns = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# --- Paragraph 1: Title ---
$p1 = $d.Paragraphs.Item(1)
$xml1 = '<w:p w:rsidR="00273E31" w:rsidRDefault="00CC4DA0"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="44"/></w:rPr><w:t>The Enigma of Life: Unveiling the Symphony of Cells and Their Remarkable Processes</w:t></w:r></w:p>'
$p1.Range.InsertXML($xml1)

# --- Paragraph 2: Author ---
$p2 = $d.Paragraphs.Item(2)
$xml2 = '<w:p w:rsidR="00273E31" w:rsidRDefault="00CC4DA0"><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>Olivia Martin</w:t></w:r></w:p>'
$p2.Range.InsertXML($xml2)

# --- Paragraph 3: Email ---
$p3 = $d.Paragraphs.Item(3)
$xml3 = '<w:p w:rsidR="00273E31" w:rsidRDefault="00CC4DA0"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>oliviamartin87@gmailnet</w:t></w:r></w:p>'
$p3.Range.InsertXML($xml3)

# --- Paragraph 5: Big body paragraph ---
$p5 = $d.Paragraphs.Item(5)
$xml5 = '<w:p w:rsidR="00273E31" w:rsidRDefault="00CC4DA0"><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>In the vast expanse of knowledge, biology stands as a captivating subject that unravels the intricate symphony of life, delving into the mysteries of living organisms and the remarkable processes that orchestrate their existence</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Biology, the science of life, is a profound exploration of the fundamental units of life, the cells, and the intricate web of interactions that govern their behaviors</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> From the depths of molecules to the heights of ecosystems, biology unveils the interplay of structures, functions, and adaptations that enable organisms to thrive in a diverse and ever-changing world</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Biology extends its reach beyond the microscopic realm, encompassing the study of organisms as individuals and as integral parts of populations and communities</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It unravels the complexities of ecosystems, examining the delicate balance between species and the intricacies of their interactions</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Biology also investigates the remarkable processes of reproduction, development, and inheritance, unveiling the mechanisms by which life perpetuates itself and evolves over time</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>The study of biology is a journey of exploration, discovery, and awe</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It invites us to ponder the tapestry of life, from the elegance of a butterfly''s wing to the resilience of a microscopic bacterium</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Through biology, we gain insights into the interconnectedness of all living things, recognizing our place in the grand symphony of life and the responsibility we bear as stewards of the planet we inhabit</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Body:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Biology unveils the intricate symphony of life through a diverse range of sub-disciplines, each contributingDu Te De Shi Jiao of this captivating subject</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Cellular biology delves into the inner workings of cells, exploring the structures, functions, and interactions of the organelles that orchestrate life''s processes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Molecular biology uncovers the secrets of DNA, RNA, and proteins, revealing the genetic code that governs heredity and guides the development and functioning of organisms</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:br/><w:t>Ecology investigates the intricate connections between organisms and their environments, examining the delicate balance of ecosystems and the interactions that shape the distribution and abundance of species</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Evolutionary biology delves into the origins of life and the remarkable journey of evolution, uncovering the mechanisms that drive adaptation, speciation, and the emergence of new forms of life</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>The study of biology has far-reaching implications for human societies</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It contributes to advancements in medicine, biotechnology, agriculture, and environmental conservation</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Understanding the fundamentals of biology empowers us to develop new treatments for diseases, enhance agricultural productivity, and devise strategies for preserving the fragile ecosystems upon which we depend</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Biology also plays a vital role in addressing global challenges such as climate change, food security, and biodiversity loss</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$p5.Range.InsertXML($xml5)

# --- Paragraph 7: Summary body ---
$p7 = $d.Paragraphs.Item(7)
$xml7 = '<w:p w:rsidR="00273E31" w:rsidRDefault="00CC4DA0"><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>Biology, the science of life, delves into the mysteries of living organisms and the intricate processes that orchestrate their existence</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> It unveils the symphony of cells, the diversity of ecosystems, and the remarkable journey of evolution</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Biology''s sub-disciplines, ranging from cellular biology to ecology, provide diverse perspectives on the captivating subject of life</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> The study of biology has profound implications for human societies, contributing to advancements in medicine, biotechnology, agriculture, and environmental conservation</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Biology empowers us to understand our place in the grand tapestry of life and equips us with the knowledge to address global challenges and strive for a sustainable future</w:t></w:r><w:r w:rsidR="00CB6CA4"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$p7.Range.InsertXML($xml7)

# --- Append trailing empty paragraph at end of body ---
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
